# Test_Cases.xlsx — "all test cases executed and updated the test_cases excel file"
#
# Fill in the "Actual Result" column (I) for the first three executed test
# cases, widen column I to fit the new text, move the selection/viewport
# over to the newly-populated area, and mark the workbook window minimized
# (mirroring the author's recorded window state on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Actual Result values for the executed test cases (rows 2-4) ---------
$ws.Range("I2").Value = "User registered successfully and confirmation notification received"
$ws.Range("I3").Value = "Error of validation is showing."
$ws.Range("I4").Value = "user logged in and notification for the login can be seen."

# --- Widen column I (Actual Result) to fit the new content ---------------
$ws.Columns.Item(9).ColumnWidth = 60.85546875

# --- Move the viewport / active selection over to the edited area -------
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("I5").Select()

# --- Reflect the saved window state (minimized) on reopen ----------------
$win.WindowState = -4140
